$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "semana 11 de 2025" update: the Poisson event table (Esperado/Observado/valor p
# per "evento") is refreshed with the new week counts. Event 330 ("Hepatitis a")
# and event 352 ("Infecciones de sitio quirurgico...") are newly added rows, and the
# old event 610 row is dropped, so every row from event 330 onward shifts down by one
# position. Simplest and most robust: rewrite every data cell (rows 2-30) directly.

# Column A holds event codes as text (e.g. "113"), not numbers - force Text format
# first so Excel does not silently coerce the numeric-looking strings to numbers.
$ws.Range("A2:A30").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = '113'
$ws.Cells.Item(2, 2).Value = 'Desnutrici”n aguda en menores de 5 anos'
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 0.37

$ws.Cells.Item(3, 1).Value = '115'
$ws.Cells.Item(3, 2).Value = 'Cancer en menores de 18 anos'
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 0.37

$ws.Cells.Item(4, 1).Value = '155'
$ws.Cells.Item(4, 2).Value = 'Cancer de la mama y cuello uterino'
$ws.Cells.Item(4, 3).Value = 4
$ws.Cells.Item(4, 4).Value = 9
$ws.Cells.Item(4, 5).Value = 0.01

$ws.Cells.Item(5, 1).Value = '210'
$ws.Cells.Item(5, 2).Value = 'Dengue'
$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 4).Value = 15
$ws.Cells.Item(5, 5).Value = 0

$ws.Cells.Item(6, 1).Value = '215'
$ws.Cells.Item(6, 2).Value = 'Defectos congenitos'
$ws.Cells.Item(6, 3).Value = 2
$ws.Cells.Item(6, 4).Value = 4
$ws.Cells.Item(6, 5).Value = 0.09

$ws.Cells.Item(7, 1).Value = '220'
$ws.Cells.Item(7, 2).Value = 'Dengue grave'
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 1

$ws.Cells.Item(8, 1).Value = '300'
$ws.Cells.Item(8, 2).Value = 'Agresiones por animales potencialmente transmisores de rabia'
$ws.Cells.Item(8, 3).Value = 46
$ws.Cells.Item(8, 4).Value = 42
$ws.Cells.Item(8, 5).Value = 0.05

$ws.Cells.Item(9, 1).Value = '330'
$ws.Cells.Item(9, 2).Value = 'Hepatitis a'
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0.37

$ws.Cells.Item(10, 1).Value = '340'
$ws.Cells.Item(10, 2).Value = 'Hepatitis b, c y coinfeccion hepatitis b y delta'
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 2
$ws.Cells.Item(10, 5).Value = 0

$ws.Cells.Item(11, 1).Value = '342'
$ws.Cells.Item(11, 2).Value = 'Enfermedades huerfanas - raras'
$ws.Cells.Item(11, 3).Value = 3
$ws.Cells.Item(11, 4).Value = 5
$ws.Cells.Item(11, 5).Value = 0.1

$ws.Cells.Item(12, 1).Value = '346'
$ws.Cells.Item(12, 2).Value = 'Ira por virus nuevo'
$ws.Cells.Item(12, 3).Value = 12
$ws.Cells.Item(12, 4).Value = 4
$ws.Cells.Item(12, 5).Value = 0.01

$ws.Cells.Item(13, 1).Value = '348'
$ws.Cells.Item(13, 2).Value = 'Infeccion respiratoria aguda grave irag inusitada'
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 1

$ws.Cells.Item(14, 1).Value = '352'
$ws.Cells.Item(14, 2).Value = 'Infecciones de sitio quirurgico asociadas a procedimiento medico quirurgico'
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0.37

$ws.Cells.Item(15, 1).Value = '355'
$ws.Cells.Item(15, 2).Value = 'Enfermedad transmitida por alimentos o agua (eta)'
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 1

$ws.Cells.Item(16, 1).Value = '356'
$ws.Cells.Item(16, 2).Value = 'Intento de suicidio'
$ws.Cells.Item(16, 3).Value = 14
$ws.Cells.Item(16, 4).Value = 14
$ws.Cells.Item(16, 5).Value = 0.11

$ws.Cells.Item(17, 1).Value = '357'
$ws.Cells.Item(17, 2).Value = 'Iad - infecciones asociadas a dispositivos - individual'
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0.14

$ws.Cells.Item(18, 1).Value = '365'
$ws.Cells.Item(18, 2).Value = 'Intoxicaciones'
$ws.Cells.Item(18, 3).Value = 8
$ws.Cells.Item(18, 4).Value = 13
$ws.Cells.Item(18, 5).Value = 0.03

$ws.Cells.Item(19, 1).Value = '455'
$ws.Cells.Item(19, 2).Value = 'Leptospirosis'
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 5).Value = 0.37

$ws.Cells.Item(20, 1).Value = '465'
$ws.Cells.Item(20, 2).Value = 'Malaria'
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 0

$ws.Cells.Item(21, 1).Value = '535'
$ws.Cells.Item(21, 2).Value = 'Meningitis bacteriana y enfermedad meningoc”cica'
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 1

$ws.Cells.Item(22, 1).Value = '549'
$ws.Cells.Item(22, 2).Value = 'Morbilidad materna extrema'
$ws.Cells.Item(22, 3).Value = 5
$ws.Cells.Item(22, 4).Value = 2
$ws.Cells.Item(22, 5).Value = 0.08

$ws.Cells.Item(23, 1).Value = '560'
$ws.Cells.Item(23, 2).Value = 'Mortalidad perinatal y neonatal tardia'
$ws.Cells.Item(23, 3).Value = 1
$ws.Cells.Item(23, 4).Value = 1
$ws.Cells.Item(23, 5).Value = 0.37

$ws.Cells.Item(24, 1).Value = '580'
$ws.Cells.Item(24, 2).Value = 'Mortalidad por dengue'
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 1

$ws.Cells.Item(25, 1).Value = '620'
$ws.Cells.Item(25, 2).Value = 'Parotiditis'
$ws.Cells.Item(25, 3).Value = 1
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0.37

$ws.Cells.Item(26, 1).Value = '740'
$ws.Cells.Item(26, 2).Value = 'Sifilis congenita'
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(26, 5).Value = 0

$ws.Cells.Item(27, 1).Value = '750'
$ws.Cells.Item(27, 2).Value = 'Sifilis gestacional'
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = 0.27

$ws.Cells.Item(28, 1).Value = '813'
$ws.Cells.Item(28, 2).Value = 'Tuberculosis'
$ws.Cells.Item(28, 3).Value = 6
$ws.Cells.Item(28, 4).Value = 5
$ws.Cells.Item(28, 5).Value = 0.16

$ws.Cells.Item(29, 1).Value = '831'
$ws.Cells.Item(29, 2).Value = 'Varicela individual'
$ws.Cells.Item(29, 3).Value = 6
$ws.Cells.Item(29, 4).Value = 7
$ws.Cells.Item(29, 5).Value = 0.14

$ws.Cells.Item(30, 1).Value = '850'
$ws.Cells.Item(30, 2).Value = 'Vih/sida/mortalidad por sida'
$ws.Cells.Item(30, 3).Value = 8
$ws.Cells.Item(30, 4).Value = 7
$ws.Cells.Item(30, 5).Value = 0.14
